# Mohammad Darmanloo's (row 8) quiz scores were recorded; previously the
# score cells F8:K8 were left blank. Fill them in now that his grade is
# known, and update the active cell selection to reflect where the user
# was working (M8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F8:K8 -> HW3 Q1, Q2, Q3, Q4, HW4 Q1, Q2 for mohammad darmanloo
$ws.Range("F8").Value = 125
$ws.Range("G8").Value = 97
$ws.Range("H8").Value = 92
$ws.Range("I8").Value = 90
$ws.Range("J8").Value = 105
$ws.Range("K8").Value = 100

$ws.Activate()
$ws.Range("M8").Select()
